$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) is treated as text so formatted numeric strings
# (e.g. trailing zeros, thousand-dot separators) are preserved exactly.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "28.546.93"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.910.78"
$ws.Range("E3").Value = "  +5.20%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "314.71"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.5205"
$ws.Range("E7").Value = "  +4.92%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3950"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.09678"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.152"
$ws.Range("E10").Value = "  +4.15%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "42.04"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.540"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "21.17"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.914.59"
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.548"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "94.69"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.00001135"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.06654"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "18.25"
$ws.Range("E20").Value = "  +6.30%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.317"
$ws.Range("E22").Value = "  +6.01%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "28.646.01"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "11.53"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.304"
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.685"
$ws.Range("E26").Value = "  +11.87%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.134.82"
$ws.Range("E27").Value = "  +5.52%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "21.25"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "158.57"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "128.55"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.110"
$ws.Range("E31").Value = "  +7.39%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1086"
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.747"
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.636"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "9.873"
$ws.Range("E35").Value = "  +10.90%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06793"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02433"
$ws.Range("E37").Value = "  +4.12%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.261"
$ws.Range("E38").Value = "  +7.50%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2235"
$ws.Range("E39").Value = "  +4.63%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "11.79"
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6472"
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "5.081"
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.192"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6104"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.752"
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "1.287"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.031"
$ws.Range("E49").Value = "  +4.80%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "125.07"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.207"
$ws.Range("E51").Value = "  +2.43%  "
